$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three more rows (26-28) following the same pattern as the existing
# rows (e.g. row 25): part number, description and price.
# Force the cells to be stored as text (matching the existing rows) rather
# than letting Excel auto-convert the numeric-looking values, then restore
# the default "Normal" style so no extra formatting is left behind.
for ($r = 26; $r -le 28; $r++) {
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 3).NumberFormat = "@"

    $ws.Cells.Item($r, 1).Value = "91697550"
    $ws.Cells.Item($r, 2).Value = "Ar condicionado Split 24000 BTUs Quente e Frio 220V Series A1 TCL"
    $ws.Cells.Item($r, 3).Value = "3,949,90"

    $ws.Cells.Item($r, 1).Style = "Normal"
    $ws.Cells.Item($r, 2).Style = "Normal"
    $ws.Cells.Item($r, 3).Style = "Normal"
}
